$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 44-55 entirely so their shared strings are purged from the
# workbook and can be rebuilt fresh, in the exact order the new data below is written.
$ws.Range("A44:N55").EntireRow.Delete()

# Row 44: Bandits Twsited Arm Port
$ws.Range("A44").Value = "Bandits Twsited Arm Port"
$ws.Range("B44").Value = "Twisted Memories"
$ws.Range("E44").Value = "A port formed by a man whos arm was twisted by dark magics performed on him by the Church of God"
$ws.Range("F44").Value = 1
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 5
$ws.Range("I44").Value = 1
$ws.Range("J44").Value = 448
$ws.Range("K44").Value = 128
$ws.Range("M44").Value = "Yes"

# Row 45: Church of God
$ws.Range("A45").Value = "Church of God"
$ws.Range("B45").Value = "Twisted Memories"
$ws.Range("E45").Value = "A church controlled by the Twsited Bishop who controls the land through persecution and supression of free will."
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 6
$ws.Range("I45").Value = 1
$ws.Range("J45").Value = 384
$ws.Range("K45").Value = 304
$ws.Range("M45").Value = "Yes"

# Row 46: Twsited grave site
$ws.Range("A46").Value = "Twsited grave site"
$ws.Range("B46").Value = "Twisted Memories"
$ws.Range("E46").Value = "A grave site of a man unknown and unamed. The land here is twisted, broken and shattered. Who lies burried here?"
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 5
$ws.Range("I46").Value = 1
$ws.Range("J46").Value = 384
$ws.Range("K46").Value = 32
$ws.Range("M46").Value = "Yes"

# Row 47: Federation Controlled Town
$ws.Range("A47").Value = "Federation Controlled Town"
$ws.Range("B47").Value = "Delusional Memories"
$ws.Range("E47").Value = "The people here are slaves to the Federation. There's the poor and the rich, there is no in between. You either work for The Federation for nothing but scraps, or you are born into a family of one who works for The Federation. Even soldiers sent to die have a higher standing then the people of this town."
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = 6
$ws.Range("I47").Value = 1
$ws.Range("J47").Value = 16
$ws.Range("K47").Value = 224
$ws.Range("M47").Value = "Yes"

# Row 48: Delusional Abandoned Gold Mines
$ws.Range("A48").Value = "Delusional Abandoned Gold Mines"
$ws.Range("B48").Value = "Delusional Memories"
$ws.Range("E48").Value = "These old Gold Mines hold the memories of the past as haunting apperations."
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 6
$ws.Range("I48").Value = 1
$ws.Range("J48").Value = 448
$ws.Range("K48").Value = 320
$ws.Range("M48").Value = "Yes"

# Row 49: Twisted Dimensional Gate
$ws.Range("A49").Value = "Twisted Dimensional Gate"
$ws.Range("B49").Value = "Hell"
$ws.Range("D49").Value = "Twisted Tree Branch"
$ws.Range("E49").Value = "A place that leads to the delisions of those who live in a fantasy world. Unable to esxape their past, unwilling to let go. The world that once was now becomes twisted and broken in the eyes of The Wondering Prince"
$ws.Range("G49").Value = 1
$ws.Range("J49").Value = 464
$ws.Range("K49").Value = 64
$ws.Range("L49").Value = 6
$ws.Range("M49").Value = "No"

# Row 50: Brothels Port
$ws.Range("A50").Value = "Brothels Port"
$ws.Range("B50").Value = "Twisted Memories"
$ws.Range("E50").Value = "A port town that is dilapidated and run down. It is controlled by the Madams who service the clerics who come through."
$ws.Range("F50").Value = 1
$ws.Range("G50").Value = 1
$ws.Range("I50").Value = 1
$ws.Range("J50").Value = 288
$ws.Range("K50").Value = 304
$ws.Range("M50").Value = "No"

# Row 51: Hermits Port
$ws.Range("A51").Value = "Hermits Port"
$ws.Range("B51").Value = "Twisted Memories"
$ws.Range("D51").Value = "Old Silver Flask"
$ws.Range("E51").Value = "Home to the hermit. The Drunk. He lives alone and calls this port his home."
$ws.Range("F51").Value = 1
$ws.Range("G51").Value = 1
$ws.Range("I51").Value = 1
$ws.Range("J51").Value = 32
$ws.Range("K51").Value = 304
$ws.Range("M51").Value = "No"

# Row 52: Northren Port
$ws.Range("A52").Value = "Northren Port"
$ws.Range("B52").Value = "Delusional Memories"
$ws.Range("E52").Value = "A port to the north. The men and women here are hardened fighters and vetrans of the battles between the Federation and free people of the south."
$ws.Range("F52").Value = 1
$ws.Range("G52").Value = 1
$ws.Range("I52").Value = 1
$ws.Range("J52").Value = 304
$ws.Range("K52").Value = 192
$ws.Range("M52").Value = "No"

# Row 53: Southren Port
$ws.Range("A53").Value = "Southren Port"
$ws.Range("B53").Value = "Delusional Memories"
$ws.Range("E53").Value = "A port of the free people in the south. They are known to trade with peoples of the Northren Port for supplies and sending goods and services back and forth. Although these days with The Federation it's harder to trade."
$ws.Range("F53").Value = 1
$ws.Range("G53").Value = 1
$ws.Range("I53").Value = 1
$ws.Range("J53").Value = 304
$ws.Range("K53").Value = 288
$ws.Range("M53").Value = "No"

# Row 54: Alchemy Corrupted Church
$ws.Range("A54").Value = "Alchemy Corrupted Church"
$ws.Range("B54").Value = "Delusional Memories"
$ws.Range("D54").Value = "Purgatory's Lantern"
$ws.Range("E54").Value = "A church corrupted by the magics of Alchemy. The gates of time have opened here, the judges of time step forth."
$ws.Range("G54").Value = 1
$ws.Range("J54").Value = 400
$ws.Range("K54").Value = 16
$ws.Range("L54").Value = 7
$ws.Range("M54").Value = "No"

# Row 55: Federation City
$ws.Range("A55").Value = "Federation City"
$ws.Range("B55").Value = "Delusional Memories"
$ws.Range("E55").Value = "The main city where the Federation organizes it's military plans from. No army, not even The Red Hawks have managed to breech the city because of the Alchemy that The Church practices, in conjunction with thier religious beliefs."
$ws.Range("G55").Value = 1
$ws.Range("J55").Value = 80
$ws.Range("K55").Value = 96
$ws.Range("M55").Value = "No"

# Row 56: Abandonded Chapel
$ws.Range("A56").Value = "Abandonded Chapel"
$ws.Range("E56").Value = "An old decrepid chapel in the middle of no where. Half burned, half rotted, what remains is a story of the past."
$ws.Range("G56").Value = 1
$ws.Range("H56").Value = 3
$ws.Range("I56").Value = 1
$ws.Range("J56").Value = 208
$ws.Range("K56").Value = 416
$ws.Range("M56").Value = "Yes"
